# Link the front-end login panel to Cognito:
# add an explanatory text box to the "Donate" slide describing the
# Login/Registration panel and its relationship to the Donation panel.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# EMU -> point conversion (1 pt = 12700 EMU) since Shapes.AddTextbox
# takes Left/Top/Width/Height in points.
$left   = 1242927 / 12700
$top    = 1259353 / 12700
$width  = 9669643 / 12700
$height = 646331 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"

$tb.Fill.Visible = 0
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tb.TextFrame.TextRange.Text = "Login/Registration Panel"
$null = $tb.TextFrame.TextRange.InsertAfter("`rDonation Panel: The user can only donate a book when they have logged in.")
